# Update "想去人数" (column F) counts on the "展览", "演出" and "全部类型" sheets
# to reflect newly generated output (commit: "Update gh-pages to output
# generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibitions) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value  = 3304
$wsExhibit.Range("F5").Value  = 2391
$wsExhibit.Range("F8").Value  = 1357
$wsExhibit.Range("F10").Value = 289
$wsExhibit.Range("F11").Value = 503
$wsExhibit.Range("F14").Value = 94
$wsExhibit.Range("F16").Value = 8340
$wsExhibit.Range("F17").Value = 363
$wsExhibit.Range("F23").Value = 574
$wsExhibit.Range("F25").Value = 1151
$wsExhibit.Range("F27").Value = 1938
$wsExhibit.Range("F28").Value = 1936
$wsExhibit.Range("F30").Value = 1720
$wsExhibit.Range("F36").Value = 72
$wsExhibit.Range("F40").Value = 221
$wsExhibit.Range("F41").Value = 394
$wsExhibit.Range("F42").Value = 73
$wsExhibit.Range("F44").Value = 248

# --- Sheet: 演出 (Performances) ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F5").Value = 14

# --- Sheet: 全部类型 (All Types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value  = 3304
$wsAll.Range("F6").Value  = 2391
$wsAll.Range("F9").Value  = 1357
$wsAll.Range("F12").Value = 289
$wsAll.Range("F13").Value = 503
$wsAll.Range("F15").Value = 94
$wsAll.Range("F17").Value = 8340
$wsAll.Range("F18").Value = 363
$wsAll.Range("F20").Value = 14
$wsAll.Range("F25").Value = 574
$wsAll.Range("F27").Value = 1151
$wsAll.Range("F29").Value = 1938
$wsAll.Range("F30").Value = 1936
$wsAll.Range("F31").Value = 1720
$wsAll.Range("F37").Value = 72
$wsAll.Range("F41").Value = 221
$wsAll.Range("F42").Value = 394
$wsAll.Range("F47").Value = 73
$wsAll.Range("F49").Value = 248
